# Weekly crime data refresh: Volume 30 Number 32 -> 33, week 8/7-8/13 -> 8/14-8/20,
# and refreshed crime-complaint figures (rows 14-30) for the new reporting week.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number and reporting week ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Crime complaints table (rows 14-30): refreshed Week/28-Day/YTD/2-Year figures ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = -12.5
$ws.Range("I14").Value = 80
$ws.Range("J14").Value = 91
$ws.Range("K14").Value = -12.087912087912
$ws.Range("L14").Value = -13.978494623655
$ws.Range("M14").Value = -6.976744186046
$ws.Range("N14").Value = -74.522292993630
$ws.Range("C15").Value = 13
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 160
$ws.Range("F15").Value = 26
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = -7.142857142857
$ws.Range("I15").Value = 247
$ws.Range("J15").Value = 257
$ws.Range("K15").Value = -3.891050583657
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = 26.666666666666
$ws.Range("N15").Value = -45.714285714285
$ws.Range("C16").Value = 112
$ws.Range("D16").Value = 104
$ws.Range("E16").Value = 7.692307692307
$ws.Range("F16").Value = 446
$ws.Range("G16").Value = 472
$ws.Range("H16").Value = -5.508474576271
$ws.Range("I16").Value = 3052
$ws.Range("J16").Value = 3261
$ws.Range("K16").Value = -6.409076970254
$ws.Range("L16").Value = 34.094903339191
$ws.Range("M16").Value = 11.468224981738
$ws.Range("N16").Value = -69.892473118279
$ws.Range("C17").Value = 164
$ws.Range("D17").Value = 125
$ws.Range("E17").Value = 31.2
$ws.Range("F17").Value = 682
$ws.Range("G17").Value = 581
$ws.Range("H17").Value = 17.383820998278
$ws.Range("I17").Value = 5127
$ws.Range("J17").Value = 4678
$ws.Range("K17").Value = 9.598118854211
$ws.Range("L17").Value = 34.109338216060
$ws.Range("M17").Value = 80.084299262381
$ws.Range("N17").Value = -12.776454576386
$ws.Range("C18").Value = 59
$ws.Range("D18").Value = 43
$ws.Range("E18").Value = 37.209302325581
$ws.Range("F18").Value = 219
$ws.Range("G18").Value = 187
$ws.Range("H18").Value = 17.112299465240
$ws.Range("I18").Value = 1898
$ws.Range("J18").Value = 1850
$ws.Range("K18").Value = 2.594594594594
$ws.Range("L18").Value = 42.385596399099
$ws.Range("M18").Value = -6.132542037586
$ws.Range("N18").Value = -84.129107784931
$ws.Range("C19").Value = 175
$ws.Range("D19").Value = 166
$ws.Range("E19").Value = 5.421686746987
$ws.Range("F19").Value = 674
$ws.Range("G19").Value = 656
$ws.Range("H19").Value = 2.743902439024
$ws.Range("I19").Value = 4958
$ws.Range("J19").Value = 5066
$ws.Range("K19").Value = -2.131859455191
$ws.Range("L19").Value = 24.043032274205
$ws.Range("M19").Value = 70.906583936573
$ws.Range("N19").Value = 5.332483535160
$ws.Range("C20").Value = 111
$ws.Range("D20").Value = 60
$ws.Range("E20").Value = 85
$ws.Range("F20").Value = 444
$ws.Range("G20").Value = 271
$ws.Range("H20").Value = 63.837638376383
$ws.Range("I20").Value = 3419
$ws.Range("J20").Value = 2501
$ws.Range("K20").Value = 36.705317872850
$ws.Range("L20").Value = 100.881316098707
$ws.Range("M20").Value = 158.232628398792
$ws.Range("N20").Value = -64.705275110973
$ws.Range("C21").Value = 636
$ws.Range("D21").Value = 504
$ws.Range("E21").Value = 26.190476190476
$ws.Range("F21").Value = 2498
$ws.Range("G21").Value = 2203
$ws.Range("H21").Value = 13.390830685429
$ws.Range("I21").Value = 18781
$ws.Range("J21").Value = 17704
$ws.Range("K21").Value = 6.083370989606
$ws.Range("L21").Value = 39.614927148379
$ws.Range("M21").Value = 55.048295220011
$ws.Range("N21").Value = -56.46197000255
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = -62.5
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = -12.5
$ws.Range("I22").Value = 181
$ws.Range("J22").Value = 226
$ws.Range("K22").Value = -19.911504424778
$ws.Range("L22").Value = 21.476510067114
$ws.Range("M22").Value = -12.135922330097
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 30
$ws.Range("E23").Value = 6.666666666666
$ws.Range("F23").Value = 138
$ws.Range("G23").Value = 116
$ws.Range("H23").Value = 18.965517241379
$ws.Range("I23").Value = 1135
$ws.Range("J23").Value = 1026
$ws.Range("K23").Value = 10.623781676413
$ws.Range("L23").Value = 50.132275132275
$ws.Range("M23").Value = 65.211062590975
$ws.Range("C24").Value = 355
$ws.Range("D24").Value = 385
$ws.Range("E24").Value = -7.792207792207
$ws.Range("F24").Value = 1404
$ws.Range("G24").Value = 1519
$ws.Range("H24").Value = -7.570770243581
$ws.Range("I24").Value = 11380
$ws.Range("J24").Value = 11734
$ws.Range("K24").Value = -3.016874041247
$ws.Range("L24").Value = 43.886711341509
$ws.Range("M24").Value = 40.824155426308
$ws.Range("C25").Value = 206
$ws.Range("D25").Value = 182
$ws.Range("E25").Value = 13.186813186813
$ws.Range("F25").Value = 846
$ws.Range("G25").Value = 769
$ws.Range("H25").Value = 10.013003901170
$ws.Range("I25").Value = 6742
$ws.Range("J25").Value = 6450
$ws.Range("K25").Value = 4.527131782945
$ws.Range("L25").Value = 26.420401275079
$ws.Range("M25").Value = -5.162470108313
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 166.666666666667
$ws.Range("F26").Value = 57
$ws.Range("G26").Value = 47
$ws.Range("H26").Value = 21.276595744680
$ws.Range("I26").Value = 421
$ws.Range("J26").Value = 444
$ws.Range("K26").Value = -5.180180180180
$ws.Range("L26").Value = 13.783783783783
$ws.Range("C27").Value = 24
$ws.Range("D27").Value = 27
$ws.Range("E27").Value = -11.111111111111
$ws.Range("F27").Value = 74
$ws.Range("G27").Value = 80
$ws.Range("H27").Value = -7.5
$ws.Range("I27").Value = 664
$ws.Range("J27").Value = 584
$ws.Range("K27").Value = 13.698630136986
$ws.Range("L27").Value = 18.571428571428
$ws.Range("C28").Value = 14
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 250
$ws.Range("F28").Value = 35
$ws.Range("G28").Value = 32
$ws.Range("H28").Value = 9.375
$ws.Range("I28").Value = 263
$ws.Range("J28").Value = 335
$ws.Range("K28").Value = -21.492537313432
$ws.Range("L28").Value = -30.789473684210
$ws.Range("M28").Value = -15.705128205128
$ws.Range("N28").Value = -70.874861572536
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 29
$ws.Range("G29").Value = 28
$ws.Range("H29").Value = 3.571428571428
$ws.Range("I29").Value = 215
$ws.Range("J29").Value = 287
$ws.Range("K29").Value = -25.087108013937
$ws.Range("L29").Value = -33.846153846153
$ws.Range("M29").Value = -16.988416988417
$ws.Range("N29").Value = -73.748473748473
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 14
$ws.Range("J30").Value = 34
$ws.Range("K30").Value = -58.823529411764
$ws.Range("L30").Value = -56.25

# C14 and C30 previously held a literal placeholder text "0"; now numeric data is available.
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C30").NumberFormat = "#,##0"
